$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add a new entry row (row 29) recording work on 21/01/2018 by Giovanni
$ws.Range("A29").Value = Get-Date -Year 2018 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("B29").Value = "Giovanni"
$ws.Range("C29").Value = "Implementati tasti gestione memoria. Implementato uso di combina. Cambiata implementazione della memoria nel modello ed alcuni metodi del controller."
$ws.Range("D29").Value = 0.25
$ws.Range("E29").Value = 0.041666666666666664

# Match formatting of the row above (row 28) for columns A:E
$ws.Range("A28:E28").Copy() | Out-Null
$ws.Range("A29:E29").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(29).RowHeight = 100.8

$ws.Range("D30").Select() | Out-Null
